# ---------------------------------------------------------------------------
# goes16_fields.xlsx update:
#  - rename Sheet1 -> "Raw Data Fields"
#  - add a new sheet "Groups and Datasets" after it, with a Group/Dataset
#    comparison table (GOES-16 vs SEVIRI_m11)
#  - tidy up sheet selections / active tab
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- rename the original sheet -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Raw Data Fields"

# --- add the new sheet right after it -------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Groups and Datasets"

# --- column widths ----------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 38.998697916666664
$ws2.Columns.Item(2).ColumnWidth = 48.830729166666664
$ws2.Columns.Item(3).ColumnWidth = 49.830729166666664

# --- table content -----------------------------------------------------------
$metaB = "datetime, height_above_mean_sea_level, latitude, longitude, record_number, scan_position, sensor_azimuth_angle, sensor_view_angle, sensor_zenith_angle, solar_azimuth_angle, solar_zenith_angle, time "

$rows = @(
  @("Group", "Datasets (SEVIRI_m11)", "Datasets ported from IODAv1 converter(GOES-16)"),
  @("Root", "nchans, ndatetime, nlocs, nstring, nvars", "nlocs, nobs, nrecs, nvars"),
  @("EffectiveError", "brightness_temperature", ""),
  @("EffectiveQC", "brightness_temperature", ""),
  @("MetaData", $metaB, "Elevation_Angle, Scan_Angle, latitude, longitude, time"),
  @("ObsBias", "brightness_temperature", ""),
  @("ObsError", "brightness_temperature", "radiance"),
  @("ObsValue", "brightness_temperature", "radiance"),
  @("PreQC", "brightness_temperature", "radiance"),
  @("VarMetaData", "ObsError, mean_lapse_rate, polarization, sensor_band_central_radiation_frequency, sensor_band_central_radiation_wavenumber, sensor_channel, variable_names", ""),
  @("constantPredictor", "brightness_temperature", ""),
  @("emissivityPredictor", "brightness_temperature", ""),
  @("lapse_ratePredictor", "brightness_temperature", ""),
  @("lapse_rate_order_2Predictor", "brightness_temperature", ""),
  @("scan_anglePredictor", "brightness_temperature", ""),
  @("scan_angle_order_2Predictor", "brightness_temperature", ""),
  @("scan_angle_order_3Predictor", "brightness_temperature", ""),
  @("scan_angle_order_4Predictor", "brightness_temperature", "")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 1
  $ws2.Cells.Item($r, 1).Value = $rows[$i][0]
  $ws2.Cells.Item($r, 2).Value = $rows[$i][1]
  if ($rows[$i][2] -ne "") {
    $ws2.Cells.Item($r, 3).Value = $rows[$i][2]
  }
}

# --- formatting --------------------------------------------------------------

# Row 1 (header): bold, centered, B1 also wraps
$hdr = $ws2.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$ws2.Range("B1").WrapText = $true

# Row 2 ("Root" sub-header): left aligned, normal weight
$sub = $ws2.Range("A2:C2")
$sub.HorizontalAlignment = -4131   # xlLeft
$ws2.Range("B2").WrapText = $true

# Data rows 3-18: column B wraps
$ws2.Range("B3:B18").WrapText = $true

# Outer/inner borders for the header row
$ws2.Range("A1:B1").Borders.Item(8).Weight = -4138   # medium top
$ws2.Range("A1").Borders.Item(7).Weight = -4138      # medium left
$ws2.Range("A1").Borders.Item(10).Weight = 2         # thin right
$ws2.Range("B1").Borders.Item(7).Weight = 2          # thin left
$ws2.Range("B1").Borders.Item(10).Weight = 2         # thin right

# Full thin box around every populated cell (A1:C18)
$ws2.Range("A1:C18").Borders.Item(7).Weight = 2
$ws2.Range("A1:C18").Borders.Item(8).Weight = 2
$ws2.Range("A1:C18").Borders.Item(9).Weight = 2
$ws2.Range("A1:C18").Borders.Item(10).Weight = 2
$ws2.Range("A1:C18").Borders.Item(11).Weight = 2
$ws2.Range("A1:C18").Borders.Item(12).Weight = 2

# Re-apply the stronger header emphasis after the generic thin grid
$ws2.Range("A1:B1").Borders.Item(8).Weight = -4138   # medium top
$ws2.Range("A1").Borders.Item(7).Weight = -4138      # medium left

# Row heights for the two wrapped-text rows
$ws2.Rows.Item(5).RowHeight = 89
$ws2.Rows.Item(10).RowHeight = 68

Write-Output "sheets ready"
